# Regenerate the "K" column (column G) values for each game row.
# The sheet previously stored a different "Strike#"-derived figure in
# column G; this regenerates the save_data so column G holds the true K
# (strikeout) value pulled per-row, matching the header label "K" in G1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (column G), recomputed for this player.
$kValues = [ordered]@{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 0
    12 = 0
    13 = 2
    14 = 0
    15 = 1
    16 = 0
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    29 = 0
    30 = 0
    31 = 2
    32 = 2
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 2
    39 = 1
    40 = 2
    41 = 1
    42 = 1
    45 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
